$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.518312692642212
$ws.Range("B1").Value = 1.761911392211914
$ws.Range("C1").Value = 1.78416109085083
$ws.Range("D1").Value = 2.19510555267334
$ws.Range("E1").Value = 3.21607780456543
